# TrialsSetup.xlsx update (2025-12-09 12:00)
# - CADANCE trial "Progress" value corrected from 12 to 62
# - HORIZON OLE trial "Progress" value filled in as 0 (was blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 62
$ws.Range("B11").Value = 0
